$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New otolith rows (Trawl 2, Maurolicus muelleri) appended below the existing
# data, which previously ended at row 19.
$rows = @{
    20 = @("34mm SL", 1)
    21 = @("27mm SL", 1)
    22 = @("28mm SL", 2)
    23 = @("31mm SL", 2)
    24 = @("31mm SL", 1)
    25 = @("29mm SL", 2)
    26 = @("35mm SL", 2)
    27 = @("30mm SL", 2)
    28 = @("28mm SL", 1)
    29 = @("33mm SL", 2)
    30 = @("33mm SL", 2)
    31 = @("26mm SL", 2)
    32 = @("29mm SL", 2)
    33 = @("31mm SL", 2)
    34 = @("29mm SL", 2)
    35 = @("34mm SL", 2)
    36 = @("31mm SL", 2)
    37 = @("29mm SL", 2)
    38 = @("34mm SL", 2)
    39 = @("42mm SL", 1)
    40 = @("30mm SL", 2)
    41 = @("32mm SL", 2)
    42 = @("32mm SL", 2)
    43 = @("28mm SL", 2)
    44 = @("28mm SL", 2)
    45 = @("33mm SL", 2)
    46 = @("32mm SL", 2)
    47 = @("31mm SL", 2)
    48 = @("31mm SL", 2)
    49 = @("28mm SL", 2)
    50 = @("NR", 2)
    51 = @("32mm SL", 2)
    52 = @("27mm SL", 2)
    53 = @("29mm SL", 2)
    54 = @("30mm SL", 2)
    55 = @("25mm SL", 2)
}

# The workbook's shared-string table records new strings in the order they were
# first typed by the original author, which is not strictly row order (e.g. the
# "35mm SL" used on row 26 was entered before the "29mm SL" used on row 25).
# Write column C (Length) for the first row of each new value in that historical
# order first, so the regenerated shared-string table matches.
$lengthIntroOrder = @(20, 21, 22, 23, 26, 25, 27, 29, 31, 39, 41, 55)
foreach ($r in $lengthIntroOrder) {
    $ws.Cells.Item($r, 3).Value = $rows[$r][0]
}

# Now fill in every cell (A, B, C, D) for all 36 new rows in normal order.
for ($r = 20; $r -le 55; $r++) {
    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = "Maurolicus muelleri"
    $ws.Cells.Item($r, 2).Font.Italic = $true
    $ws.Cells.Item($r, 3).Value = $rows[$r][0]
    $ws.Cells.Item($r, 4).Value = $rows[$r][1]
}

# Restore the view/selection state recorded in the edited workbook: scrolled so
# row 19 is at the top, with D19 selected.
$excel.Goto($ws.Range("A19"), $true)
$ws.Range("D19").Select() | Out-Null
